$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 168, pushing existing rows 168..270 down to 169..271
$ws.Rows.Item(168).Insert()

# Populate the newly inserted row 168 with the new record's data
$ws.Range("A168").Value = 4
$ws.Range("B168").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C168").Value = "Los Lagos"
$ws.Range("D168").Value = 44777
$ws.Range("D168").NumberFormat = $ws.Range("D169").NumberFormat
$ws.Range("E168").Value = 10
$ws.Range("F168").Value = 100112032
$ws.Range("G168").Value = "Zapallo italiano"
$ws.Range("H168").Value = "Sin especificar"
$ws.Range("I168").Value = "Primera"
$ws.Range("J168").Value = 140
$ws.Range("K168").Value = 24000
$ws.Range("L168").Value = 25000
$ws.Range("M168").Value = 24500
$ws.Range("N168").Value = "$/caja 50 unidades"
$ws.Range("O168").Value = "Región de Arica y Parinacota"
$ws.Range("P168").Value = 490
$ws.Range("Q168").Value = 50
$ws.Range("R168").Value = "Hortaliza"
